$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 30 data (date 2025-08-05), appended after existing row 29
$ws.Range("A30").Value = 45874
$ws.Range("A30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B30").Value = 6299.18994140625
$ws.Range("C30").Value = 6346
$ws.Range("D30").Value = 6289.3701171875
$ws.Range("E30").Value = 6336.6298828125
$ws.Range("F30").Value = 5517410000
$ws.Range("G30").Value = -0.0048578659963033
$ws.Range("H30").Value = -1
$ws.Range("I30").Value = 'tpg re finance trust q eps beat prologis see narrow window before warehouse rent increase ap business summarybrief at pm edt tesla award ceo musk million of share valued at about billion fox nation greenlights two new season of god family football tesla award ceo musk million of share valued at about billion ap news nyse martin say ipo market open for pretty much all sector exxonmobil q eps beat but stock fell magnera announces participation at the ubs global material conference magnera announces participation at the ubs global material conference danville cvb award in tourism grant panama port challenge could turn into trump win man accused in atm robbery bombing arrested in newport beach wipfli secures strategic investment from new mountain capital wild brawl erupts between ton of parent during youth football game all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend poltreg dos first patient in phase ii trial for presymptomatic type diabetes with ptg what a september rate cut would mean for stock plus a bullish call on ai spending clash of trade reality show aim to boost prestige of u manufacturing why gen z driver suffer the highest rate of crash injury why gen z driver suffer the highest rate of crash injury why gen z driver suffer the highest rate of crash injury why gen z driver suffer the highest rate of crash injury why gen z driver suffer the highest rate of crash injury palantir earnings due a nvidia broadcom shake off trump tariff investor business daily bitcoin news today m outflow hit digital asset amid fed rate fear fanzone scoreboard the ultimate gameday upgrade for nfl fan amazon metal water hose is absolutely a life changer apple io why you should update your iphone today naphtha cushion china q oil demand decline oxford bank corporation announces six month and second quarter operating result fanzone scoreboard the ultimate gameday upgrade for nfl fan thousand of boeing worker who build fighter jet and weapon go on strike palantirs ultraexpensive valuation spark worry into result yahoo finance agibank raise brl billion in debenture offering to strengthen it funding strategy and support sustainable growth kremlin urge careful use of nuclear rhetoric a medvedev pivot to saying trump blackmailed by mossad recent home sale in the area bond yield continue to decline to lowest level since april healthcare workforce management market surge amid drive for operational efficiency and cost reduction tx rail product inc report financial result for third quarter of fiscal delta say ai speed up market analysis to price airfare letter potential for a repeat firefly aerospace lift ipo range that would value company at over billion cftc issue noaction relief for sefs on order book obligation sarat sethi inflation not coming down ha put fed in quandary nyse texas tap former abbott official a president shareholder investigation halper sadeh llc investigates sptn base rkda on behalf of shareholder how is the trump bill changing plan bnb integration announced for q in apirones cryptogateway search for suspect in tennessee killing focus near university source say a rd arrest announced hhi water quality laboratory receives excellence recognition shareholder investigation halper sadeh llc investigates sptn base rkda on behalf of shareholder all major la vega strip casino are unionized defying national trend homestyle direct medically tailored meal delivery score high in pennsylvania all major la vega strip casino are unionized defying national trend homestyle direct expands access to medically tailored meal in california hireright appoints jim desmond a chief information security officer what may happen to stock market next pridestaff launch new fort worth west office to strengthen workforce solution in north texas live will palantir post a big earnings beat after the bell local husbandandwife team to lead pridestaff asheville office the ultimate guide to top cryptos to join now for massive roi ww grainger analyst cut their forecast after q earnings skylight health launch valuebased multispecialty platform mark herr communication issue statement regarding the notice of appeal by anne pramaggiore mark herr communication issue statement regarding the notice of appeal by anne pramaggiore stock that soared the last time market were this wobbly tesla gain on elon musk new pay package is tsla stock a buy euro strengthens a postpayroll fallout set stage for fed september rate cut legendary investor vinod khosla advises gen z to invest in this one skill because chatgpt can teach you everything else fortune roadwork to alter traffic on ssu campus these analyst raise their forecast on nvent electric after strong q earnings ap business summarybrief at pm edt wall street rally and u stock recover much of friday wipeout wall street rally and u stock recover much of friday wipeout wall street rally and u stock recover much of friday wipeout wall street rally and u stock recover much of friday wipeout microbiome sequencing market on the rise cagr forecast why opendoor technology stock wa soaring again today trump will raise india tariff rate substantially over russian oil purchasesheres why india might not budge these analyst increase their forecast on regeneron pharmaceutical after betterthanexpected earnings fed rate cut now likely a soft data erodes confidence in economic resilience elon musk unveils grok update after controversial post by the ai chatbot why nvidia is the biggest risk to the stock marketnot tariff barrons dude perfect partner with regal cinema to launch dude perfect the hero tour movie in hundred of theater worldwide sp could plummet say morgan stanley tipranks cal thomas starvation in gaza whom to blame former football coach enters georgia senate race european union say it pausing retaliatory tariff against the u amid trade talk gerken named sterling professor of law dhs say daca recipient are illegal alien tell them to selfdeport ap news summary at pm edt cbo republican megabill to cost t due to higher borrowing cost state department may require visa applicant to post bond of up to to enter the u state department may require visa applicant to post bond of up to to enter the u state department may require visa applicant to post bond of up to to enter the u all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend mitch mcconnells legacy come under fire in kentucky race to replace him in the senate davenport open application for round of dream home improvement funding davenport open application for round of dream home improvement funding trump say sweeneys ad is hottest american eagle stock go wild a the great new jersey denizen yogi berra observed it deja vu all over again noncompete clause are under fire in new jersey once more how many woman on the board fewer company are saying mtg sour on republican party claim gop turned it back on america first for tax paid minnesota offer a taxpayer receipt for where dollar go pic of woman at tx capitol fighting to use womens bathroom sum up how dumb transmovement really is search for suspect in tennessee killing focus near university source say a rd arrest announced all major la vega strip casino are unionized defying national trend all major la vega strip casino are unionized defying national trend sex offender harassed marjorie taylor greene turn on gop for abandoning maga base officer find teen girl when arresting man for sex trafficking mpd say mamdani study america worst mayor brandon johnson to avoid his political pitfall report now i love her ad trump praise sydney sweeney amid jean ad fallout ap news summary at pm edt ap news summary at pm edt donald trump and the deconstruction of america mitch mcconnells legacy come under fire in kentucky race to replace him in the senate men accused of shooting up woman glenburn home and car trump shootthemessenger tactic will only hurt the economy more america best idea is under fierce attack u rep kelly morrison introduces legislation to boost police recruitment funding these epstein revelation show the deep state is real and it work for trump mary moriarty dont accept simplistic narrative on carjacking prosecution texas dems flee to illinois to block redistricting vote texas dems flee to illinois to block redistricting vote copy famine in gaza ocean city mayor vow further effort on ebike safety after nearly hitting rider with democrat out of state whats next for texas special session trump threatens to substantially raise tariff on indian good a it continues to buy russian oil trump threatens to substantially raise tariff on indian good a it continues to buy russian oil mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate every minute is crucial here israeli journalist on hostage held inside gaza michta published in realcleardefense and featured in realclearworld on nato and european security nyt claim many jewish voter back mamdani you wont believe who they interviewed college grad compete with ai for entrylevel job guest commentary cpuc dont saddle monterey peninsula ratepayer with desalination project celtic minority owner offer recordbreaking amount to buy move wnba team report mr william bill sander dubose sr how the supreme court kickstarted the redistricting arm race minnesota department of correction expands online application for visitor at adult facility state department may require visa applicant to post bond of up to to enter the u trump administration wont say why it transferred ghislaine maxwell to a minimumsecurity prison new law to make new hampshire library record available to parent grammywinner doechii bringing tour to san diego mortgage lender stick with bidenera reappraisal guidance white south african farm owner on trial for allegedly murdering two black woman and feeding their body to pig trump to substantially raise tariff on india for buying russian oil banana republic stuff top economist react to trump firing bls statistic chief ghislaine maxwell prison transfer add to trump epstein morass look on cnn host face a lee zeldin politely take her climate change b apart is priceless watch texas hold em gov abbott threatens legal action for absent dems mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate mitch mcconnells legacy come under fire in kentucky race to replace him in the senate letter ballroom location not a good idea ap news summary at pm edt new u visa rule could require bond from visitor it feel less like advertising and more like brainwashing mainer injured when their utv rolled over in nh gravel pit former ct chief state attorney fined in ethic case what to know texas house issue arrest warrant for rogue democrat who fled to blue state rep nancy mace is latest republican to announce bid for south carolina governor minnesota congresswoman call for more police recruitment funding why dozen of democrat left texas and how republican want to punish them why dozen of democrat left texas and how republican want to punish them why dozen of democrat left texas and how republican want to punish them why dozen of democrat left texas and how republican want to punish them'
$ws.Range("J30").Value = 'tpg finance trust q eps beat prologis narrow window warehouse rent increase business pm tesla award ceo musk million valued billion fox nation greenlights season god family football tesla award ceo musk million valued billion nyse martin ipo market open pretty sector exxonmobil q eps beat stock fell magnera announces participation ubs global material conference magnera announces participation ubs global material conference danville cvb award tourism grant panama port challenge turn trump man accused atm robbery bombing arrested newport beach wipfli secures strategic investment mountain wild brawl erupts ton parent youth football game major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend poltreg do patient phase ii trial presymptomatic type diabetes ptg september rate cut stock plus bullish call ai spending clash trade reality aim boost prestige manufacturing gen z driver suffer highest rate crash injury gen z driver suffer highest rate crash injury gen z driver suffer highest rate crash injury gen z driver suffer highest rate crash injury gen z driver suffer highest rate crash injury palantir earnings due nvidia broadcom shake trump tariff investor business bitcoin outflow hit digital asset amid fed rate fear fanzone scoreboard ultimate gameday upgrade nfl fan amazon metal water hose absolutely life changer apple io update iphone naphtha cushion china q oil demand decline oxford bank corporation announces six month quarter operating fanzone scoreboard ultimate gameday upgrade nfl fan thousand boeing worker build fighter jet weapon strike palantirs ultraexpensive valuation spark worry yahoo finance agibank raise brl billion debenture offering strengthen funding strategy support sustainable growth kremlin urge careful nuclear rhetoric medvedev pivot saying trump blackmailed mossad recent sale area bond yield continue decline lowest level since april healthcare workforce market surge amid drive operational efficiency cost reduction tx rail product report financial third quarter fiscal delta ai speed market analysis price airfare potential repeat firefly aerospace lift ipo range value company billion cftc issue noaction relief sefs order book obligation sarat sethi inflation coming put fed quandary nyse texas tap abbott president shareholder investigation halper sadeh investigates sptn base rkda behalf shareholder trump bill changing plan bnb integration announced q apirones cryptogateway search suspect tennessee killing focus near university source rd arrest announced hhi water quality laboratory receives excellence recognition shareholder investigation halper sadeh investigates sptn base rkda behalf shareholder major la vega strip casino unionized defying national trend homestyle direct medically tailored meal delivery score high pennsylvania major la vega strip casino unionized defying national trend homestyle direct expands access medically tailored meal california hireright appoints jim desmond chief information security officer happen stock market next pridestaff launch fort worth west office strengthen workforce solution north texas palantir post earnings beat bell local husbandandwife team lead pridestaff asheville office ultimate guide cryptos join massive roi ww grainger analyst cut forecast q earnings skylight health launch valuebased multispecialty platform mark herr communication issue statement regarding notice appeal anne pramaggiore mark herr communication issue statement regarding notice appeal anne pramaggiore stock soared market wobbly tesla gain elon musk pay package tsla stock buy euro strengthens postpayroll fallout set stage fed september rate cut legendary investor vinod khosla advises gen z invest skill chatgpt teach everything else fortune roadwork alter traffic ssu campus analyst raise forecast nvent electric strong q earnings business pm wall street rally stock recover friday wipeout wall street rally stock recover friday wipeout wall street rally stock recover friday wipeout wall street rally stock recover friday wipeout microbiome sequencing market rise cagr forecast opendoor technology stock soaring trump raise india tariff rate substantially russian oil purchasesheres india budge analyst increase forecast regeneron pharmaceutical betterthanexpected earnings fed rate cut likely soft data erodes confidence economic resilience elon musk unveils grok update controversial post ai chatbot nvidia biggest risk stock marketnot tariff barrons dude perfect partner regal cinema launch dude perfect hero tour movie hundred theater worldwide sp plummet morgan stanley tipranks cal thomas starvation gaza blame football coach enters georgia senate race european union pausing retaliatory tariff amid trade gerken named sterling professor dhs daca recipient illegal alien selfdeport summary pm cbo republican megabill cost due higher borrowing cost department require visa applicant post bond enter department require visa applicant post bond enter department require visa applicant post bond enter major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend mitch mcconnells legacy come fire kentucky race replace senate davenport open application round dream improvement funding davenport open application round dream improvement funding trump sweeneys ad hottest american eagle stock wild great jersey denizen yogi berra observed deja vu noncompete clause fire jersey woman board fewer company saying mtg sour republican party claim gop turned america tax paid minnesota offer taxpayer receipt dollar pic woman tx capitol fighting womens bathroom sum dumb transmovement search suspect tennessee killing focus near university source rd arrest announced major la vega strip casino unionized defying national trend major la vega strip casino unionized defying national trend sex offender harassed marjorie taylor greene turn gop abandoning maga base officer find teen girl arresting man sex trafficking mpd mamdani study america worst mayor brandon johnson avoid political pitfall report love ad trump praise sydney sweeney amid jean ad fallout summary pm summary pm donald trump deconstruction america mitch mcconnells legacy come fire kentucky race replace senate men accused shooting woman glenburn car trump shootthemessenger tactic hurt economy america best idea fierce attack rep kelly morrison introduces legislation boost police recruitment funding epstein revelation deep real work trump mary moriarty dont accept simplistic narrative carjacking prosecution texas dems flee illinois block redistricting vote texas dems flee illinois block redistricting vote copy famine gaza ocean mayor vow effort ebike safety nearly hitting rider democrat whats next texas special session trump threatens substantially raise tariff indian good continues buy russian oil trump threatens substantially raise tariff indian good continues buy russian oil mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate crucial israeli journalist hostage inside gaza michta published realcleardefense featured realclearworld nato european security nyt claim jewish voter mamdani wont believe interviewed college grad compete ai entrylevel job guest commentary cpuc dont saddle monterey peninsula ratepayer desalination project celtic minority owner offer recordbreaking amount buy move wnba team report mr william bill sander dubose sr supreme kickstarted redistricting arm race minnesota department correction expands online application visitor adult facility department require visa applicant post bond enter trump administration wont transferred ghislaine maxwell minimumsecurity prison hampshire library record available parent grammywinner doechii bringing tour san diego mortgage lender stick bidenera reappraisal guidance white south african owner trial allegedly murdering black woman feeding body pig trump substantially raise tariff india buying russian oil banana republic stuff economist react trump firing bls statistic chief ghislaine maxwell prison transfer add trump epstein morass look cnn host face lee zeldin politely take climate change b apart priceless watch texas em gov abbott threatens legal action absent dems mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate mitch mcconnells legacy come fire kentucky race replace senate ballroom location good idea summary pm visa rule require bond visitor feel less advertising brainwashing mainer injured utv rolled nh gravel pit ct chief attorney fined ethic texas house issue arrest warrant rogue democrat fled blue rep nancy mace latest republican announce bid south carolina governor minnesota congresswoman call police recruitment funding dozen democrat left texas republican punish dozen democrat left texas republican punish dozen democrat left texas republican punish dozen democrat left texas republican punish'

Write-Output "row30 written"
